$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 728
$ws.Range("F4").Value = 54
$ws.Range("F5").Value = 48
$ws.Range("F6").Value = 3023
$ws.Range("F7").Value = 1690
$ws.Range("F8").Value = 2031
$ws.Range("F9").Value = 324
$ws.Range("F10").Value = 301
$ws.Range("F11").Value = 873
$ws.Range("F12").Value = 962
$ws.Range("F13").Value = 207
$ws.Range("F14").Value = 429
$ws.Range("F15").Value = 1138
$ws.Range("F17").Value = 69
$ws.Range("F19").Value = 7328
$ws.Range("F20").Value = 291
$ws.Range("F21").Value = 2062
$ws.Range("F22").Value = 194
$ws.Range("F23").Value = 202
$ws.Range("F26").Value = 507
$ws.Range("F28").Value = 1125
$ws.Range("F29").Value = 949
$ws.Range("F31").Value = 125
$ws.Range("F32").Value = 239
$ws.Range("F36").Value = 22
$ws.Range("F37").Value = 154
$ws.Range("F38").Value = 258
$ws.Range("F39").Value = 34
$ws.Range("F40").Value = 154
$ws.Range("F41").Value = 289
$ws.Range("F43").Value = 196

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 2

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 728
$ws.Range("F4").Value = 54
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 48
$ws.Range("F9").Value = 3023
$ws.Range("F10").Value = 1690
$ws.Range("F11").Value = 2031
$ws.Range("F12").Value = 324
$ws.Range("F13").Value = 301
$ws.Range("F14").Value = 873
$ws.Range("F16").Value = 962
$ws.Range("F17").Value = 207
$ws.Range("F18").Value = 429
$ws.Range("F19").Value = 1138
$ws.Range("F21").Value = 69
$ws.Range("F23").Value = 7329
$ws.Range("F24").Value = 291
$ws.Range("F25").Value = 2063
$ws.Range("F27").Value = 194
$ws.Range("F28").Value = 202
$ws.Range("F31").Value = 507
$ws.Range("F33").Value = 1125
$ws.Range("F34").Value = 949
$ws.Range("F36").Value = 125
$ws.Range("F40").Value = 22
$ws.Range("F41").Value = 154
$ws.Range("F42").Value = 258
$ws.Range("F43").Value = 34
$ws.Range("F44").Value = 154
$ws.Range("F45").Value = 289
$ws.Range("F49").Value = 196

# --- Special case: G column becomes "已售罄" (sold out) text ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G9").Value = "已售罄"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G12").Value = "已售罄"
